# Auto-generated market-data refresh for the Excalibur_Profits workbook.
# For each affected Leve row, the price/profit columns (H:N) are updated
# to reflect freshly pulled market-board data. Columns whose source value
# became unavailable (no valid price) are cleared rather than zeroed, to
# match how the scheduled runner omits unattainable profit figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder
$ws.Cells.Item(6, 8).Value = 59.1  # H6 currentAveragePrice
$ws.Cells.Item(6, 9).Value = 59.1  # I6 currentAveragePriceNQ
$ws.Cells.Item(6, 10).Value = 0  # J6 currentAveragePriceHQ
$ws.Cells.Item(6, 11).Value = 177.3  # K6 LevePriceNQ
$ws.Cells.Item(6, 12).Value = 0  # L6 LevePriceHQ
$ws.Cells.Item(6, 13).Value = -65.30000000000001  # M6 LeveProfitNQ
$ws.Cells.Item(6, 14).Value = ""  # N6 LeveProfitHQ (cleared)

# Row 33: Glazed and Confused
$ws.Cells.Item(33, 8).Value = 335.3  # H33 currentAveragePrice
$ws.Cells.Item(33, 9).Value = 250.4375  # I33 currentAveragePriceNQ
$ws.Cells.Item(33, 10).Value = 674.75  # J33 currentAveragePriceHQ
$ws.Cells.Item(33, 11).Value = 250.4375  # K33 LevePriceNQ
$ws.Cells.Item(33, 12).Value = 674.75  # L33 LevePriceHQ
$ws.Cells.Item(33, 13).Value = -21.4375  # M33 LeveProfitNQ
$ws.Cells.Item(33, 14).Value = -1132.75  # N33 LeveProfitHQ

# Row 39: Riches' Brew
$ws.Cells.Item(39, 8).Value = 699.4167  # H39 currentAveragePrice
$ws.Cells.Item(39, 9).Value = 54.88889  # I39 currentAveragePriceNQ
$ws.Cells.Item(39, 11).Value = 164.66667  # K39 LevePriceNQ
$ws.Cells.Item(39, 13).Value = 131.33333  # M39 LeveProfitNQ

# Row 42: Eye of the Beholder
$ws.Cells.Item(42, 8).Value = 273.5  # H42 currentAveragePrice
$ws.Cells.Item(42, 9).Value = 275  # I42 currentAveragePriceNQ
$ws.Cells.Item(42, 11).Value = 825  # K42 LevePriceNQ
$ws.Cells.Item(42, 13).Value = -595  # M42 LeveProfitNQ

# Row 43: Growing Is Knowing
$ws.Cells.Item(43, 8).Value = 2409  # H43 currentAveragePrice
$ws.Cells.Item(43, 10).Value = 2499.889  # J43 currentAveragePriceHQ
$ws.Cells.Item(43, 12).Value = 2499.889  # L43 LevePriceHQ
$ws.Cells.Item(43, 14).Value = -2637.889  # N43 LeveProfitHQ

# Row 48: The Sting of Conscience
$ws.Cells.Item(48, 8).Value = 1593  # H48 currentAveragePrice
$ws.Cells.Item(48, 10).Value = 1889.5  # J48 currentAveragePriceHQ
$ws.Cells.Item(48, 12).Value = 5668.5  # L48 LevePriceHQ
$ws.Cells.Item(48, 14).Value = -6252.5  # N48 LeveProfitHQ

# Row 56: Sleepless in Silvertear
$ws.Cells.Item(56, 8).Value = 1593  # H56 currentAveragePrice
$ws.Cells.Item(56, 10).Value = 1889.5  # J56 currentAveragePriceHQ
$ws.Cells.Item(56, 12).Value = 5668.5  # L56 LevePriceHQ
$ws.Cells.Item(56, 14).Value = -6736.5  # N56 LeveProfitHQ

# Row 59: Shut Up Already
$ws.Cells.Item(59, 8).Value = 299.66666  # H59 currentAveragePrice
$ws.Cells.Item(59, 10).Value = 0  # J59 currentAveragePriceHQ
$ws.Cells.Item(59, 12).Value = 0  # L59 LevePriceHQ
$ws.Cells.Item(59, 14).Value = ""  # N59 LeveProfitHQ (cleared)

# Row 61: Not Taking No for an Answer
$ws.Cells.Item(61, 8).Value = 0  # H61 currentAveragePrice
$ws.Cells.Item(61, 9).Value = 0  # I61 currentAveragePriceNQ
$ws.Cells.Item(61, 11).Value = 0  # K61 LevePriceNQ
$ws.Cells.Item(61, 13).Value = ""  # M61 LeveProfitNQ (cleared)

# Row 64: Forged from the Void
$ws.Cells.Item(64, 8).Value = 0  # H64 currentAveragePrice
$ws.Cells.Item(64, 9).Value = 0  # I64 currentAveragePriceNQ
$ws.Cells.Item(64, 10).Value = 0  # J64 currentAveragePriceHQ
$ws.Cells.Item(64, 11).Value = 0  # K64 LevePriceNQ
$ws.Cells.Item(64, 12).Value = 0  # L64 LevePriceHQ
$ws.Cells.Item(64, 13).Value = ""  # M64 LeveProfitNQ (cleared)
$ws.Cells.Item(64, 14).Value = ""  # N64 LeveProfitHQ (cleared)

# Row 67: Dodging the Draft (L)
$ws.Cells.Item(67, 8).Value = 0  # H67 currentAveragePrice
$ws.Cells.Item(67, 9).Value = 0  # I67 currentAveragePriceNQ
$ws.Cells.Item(67, 10).Value = 0  # J67 currentAveragePriceHQ
$ws.Cells.Item(67, 11).Value = 0  # K67 LevePriceNQ
$ws.Cells.Item(67, 12).Value = 0  # L67 LevePriceHQ
$ws.Cells.Item(67, 13).Value = ""  # M67 LeveProfitNQ (cleared)
$ws.Cells.Item(67, 14).Value = ""  # N67 LeveProfitHQ (cleared)

# Row 80: Cleansing the Wicked Humours
$ws.Cells.Item(80, 8).Value = 4310944.5  # H80 currentAveragePrice
$ws.Cells.Item(80, 9).Value = 8929013  # I80 currentAveragePriceNQ
$ws.Cells.Item(80, 11).Value = 26787039  # K80 LevePriceNQ
$ws.Cells.Item(80, 13).Value = -26786041  # M80 LeveProfitNQ

# Row 83: Washing Away the Sins (L)
$ws.Cells.Item(83, 8).Value = 4310944.5  # H83 currentAveragePrice
$ws.Cells.Item(83, 9).Value = 8929013  # I83 currentAveragePriceNQ
$ws.Cells.Item(83, 11).Value = 80361117  # K83 LevePriceNQ
$ws.Cells.Item(83, 13).Value = -80356125  # M83 LeveProfitNQ

# Row 97: Materia Worth
$ws.Cells.Item(97, 8).Value = 491.33334  # H97 currentAveragePrice
$ws.Cells.Item(97, 10).Value = 537  # J97 currentAveragePriceHQ
$ws.Cells.Item(97, 12).Value = 1611  # L97 LevePriceHQ
$ws.Cells.Item(97, 14).Value = -2603  # N97 LeveProfitHQ

# Row 116: Growing Up
$ws.Cells.Item(116, 8).Value = 52699.082  # H116 currentAveragePrice
$ws.Cells.Item(116, 9).Value = 81227  # I116 currentAveragePriceNQ
$ws.Cells.Item(116, 10).Value = 12760  # J116 currentAveragePriceHQ
$ws.Cells.Item(116, 11).Value = 81227  # K116 LevePriceNQ
$ws.Cells.Item(116, 12).Value = 12760  # L116 LevePriceHQ
$ws.Cells.Item(116, 13).Value = -77785  # M116 LeveProfitNQ
$ws.Cells.Item(116, 14).Value = -19644  # N116 LeveProfitHQ

# Row 125: Body over Mind
$ws.Cells.Item(125, 8).Value = 617  # H125 currentAveragePrice
$ws.Cells.Item(125, 9).Value = 414  # I125 currentAveragePriceNQ
$ws.Cells.Item(125, 11).Value = 3726  # K125 LevePriceNQ
$ws.Cells.Item(125, 13).Value = -1266  # M125 LeveProfitNQ

# Row 131: Mindful Study
$ws.Cells.Item(131, 8).Value = 4669.25  # H131 currentAveragePrice
$ws.Cells.Item(131, 9).Value = 4392.3335  # I131 currentAveragePriceNQ
$ws.Cells.Item(131, 11).Value = 13177.0005  # K131 LevePriceNQ
$ws.Cells.Item(131, 13).Value = -8137.000499999998  # M131 LeveProfitNQ

# Row 132: Fast-forwarding Flora
$ws.Cells.Item(132, 8).Value = 45922.793  # H132 currentAveragePrice
$ws.Cells.Item(132, 9).Value = 48814.53  # I132 currentAveragePriceNQ
$ws.Cells.Item(132, 10).Value = 10499  # J132 currentAveragePriceHQ
$ws.Cells.Item(132, 11).Value = 146443.59  # K132 LevePriceNQ
$ws.Cells.Item(132, 12).Value = 31497  # L132 LevePriceHQ
$ws.Cells.Item(132, 13).Value = -143913.59  # M132 LeveProfitNQ
$ws.Cells.Item(132, 14).Value = -36557  # N132 LeveProfitHQ

# Row 135: For Tired Minds
$ws.Cells.Item(135, 8).Value = 425.29166  # H135 currentAveragePrice
$ws.Cells.Item(135, 9).Value = 429.73914  # I135 currentAveragePriceNQ
$ws.Cells.Item(135, 11).Value = 3867.65226  # K135 LevePriceNQ
$ws.Cells.Item(135, 13).Value = -1332.65226  # M135 LeveProfitNQ

# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 4895.7036  # H137 currentAveragePrice
$ws.Cells.Item(137, 9).Value = 4366  # I137 currentAveragePriceNQ
$ws.Cells.Item(137, 11).Value = 13098  # K137 LevePriceNQ
$ws.Cells.Item(137, 13).Value = -10548  # M137 LeveProfitNQ

# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 2743.4238  # H138 currentAveragePrice
$ws.Cells.Item(138, 9).Value = 1110.3334  # I138 currentAveragePriceNQ
$ws.Cells.Item(138, 10).Value = 4432.8276  # J138 currentAveragePriceHQ
$ws.Cells.Item(138, 11).Value = 3331.0002  # K138 LevePriceNQ
$ws.Cells.Item(138, 12).Value = 13298.4828  # L138 LevePriceHQ
$ws.Cells.Item(138, 13).Value = 1808.9998  # M138 LeveProfitNQ
$ws.Cells.Item(138, 14).Value = -23578.4828  # N138 LeveProfitHQ

# Row 141: Remedy for Reason
$ws.Cells.Item(141, 8).Value = 1142.5555  # H141 currentAveragePrice
$ws.Cells.Item(141, 9).Value = 1154  # I141 currentAveragePriceNQ
$ws.Cells.Item(141, 10).Value = 1102.5  # J141 currentAveragePriceHQ
$ws.Cells.Item(141, 11).Value = 3462  # K141 LevePriceNQ
$ws.Cells.Item(141, 12).Value = 3307.5  # L141 LevePriceHQ
$ws.Cells.Item(141, 13).Value = 1718  # M141 LeveProfitNQ
$ws.Cells.Item(141, 14).Value = -13667.5  # N141 LeveProfitHQ

$ws = $wb.Worksheets.Item("ARM")
# Row 28: 246 Kinds of Cheese
$ws.Cells.Item(28, 8).Value = 2735.5  # H28 currentAveragePrice
$ws.Cells.Item(28, 9).Value = 2735.5  # I28 currentAveragePriceNQ
$ws.Cells.Item(28, 10).Value = 0  # J28 currentAveragePriceHQ
$ws.Cells.Item(28, 11).Value = 2735.5  # K28 LevePriceNQ
$ws.Cells.Item(28, 12).Value = 0  # L28 LevePriceHQ
$ws.Cells.Item(28, 13).Value = -2543.5  # M28 LeveProfitNQ
$ws.Cells.Item(28, 14).Value = ""  # N28 LeveProfitHQ (cleared)

# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 6580758.5  # H32 currentAveragePrice
$ws.Cells.Item(32, 9).Value = 7937327  # I32 currentAveragePriceNQ
$ws.Cells.Item(32, 11).Value = 7937327  # K32 LevePriceNQ
$ws.Cells.Item(32, 13).Value = -7937040  # M32 LeveProfitNQ

# Row 61: Dealing with the Tough Stuff
$ws.Cells.Item(61, 8).Value = 2886.2083  # H61 currentAveragePrice
$ws.Cells.Item(61, 9).Value = 2697.4167  # I61 currentAveragePriceNQ
$ws.Cells.Item(61, 10).Value = 3452.5833  # J61 currentAveragePriceHQ
$ws.Cells.Item(61, 11).Value = 2697.4167  # K61 LevePriceNQ
$ws.Cells.Item(61, 12).Value = 3452.5833  # L61 LevePriceHQ
$ws.Cells.Item(61, 13).Value = -2485.4167  # M61 LeveProfitNQ
$ws.Cells.Item(61, 14).Value = -3876.5833  # N61 LeveProfitHQ

# Row 74: As the Bolt Flies
$ws.Cells.Item(74, 8).Value = 2238.907  # H74 currentAveragePrice
$ws.Cells.Item(74, 9).Value = 1094.875  # I74 currentAveragePriceNQ
$ws.Cells.Item(74, 11).Value = 1094.875  # K74 LevePriceNQ
$ws.Cells.Item(74, 13).Value = -220.875  # M74 LeveProfitNQ

# Row 77: Heavy Metal Banned (L)
$ws.Cells.Item(77, 8).Value = 2238.907  # H77 currentAveragePrice
$ws.Cells.Item(77, 9).Value = 1094.875  # I77 currentAveragePriceNQ
$ws.Cells.Item(77, 11).Value = 5474.375  # K77 LevePriceNQ
$ws.Cells.Item(77, 13).Value = -1106.375  # M77 LeveProfitNQ

# Row 99: Home Cooking
$ws.Cells.Item(99, 8).Value = 2735.5  # H99 currentAveragePrice
$ws.Cells.Item(99, 9).Value = 2735.5  # I99 currentAveragePriceNQ
$ws.Cells.Item(99, 10).Value = 0  # J99 currentAveragePriceHQ
$ws.Cells.Item(99, 11).Value = 2735.5  # K99 LevePriceNQ
$ws.Cells.Item(99, 12).Value = 0  # L99 LevePriceHQ
$ws.Cells.Item(99, 13).Value = 259.5  # M99 LeveProfitNQ
$ws.Cells.Item(99, 14).Value = ""  # N99 LeveProfitHQ (cleared)

# Row 132: Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 2146.9688  # H132 currentAveragePrice
$ws.Cells.Item(132, 9).Value = 1925.4912  # I132 currentAveragePriceNQ
$ws.Cells.Item(132, 10).Value = 3950.4285  # J132 currentAveragePriceHQ
$ws.Cells.Item(132, 11).Value = 5776.473599999999  # K132 LevePriceNQ
$ws.Cells.Item(132, 12).Value = 11851.2855  # L132 LevePriceHQ
$ws.Cells.Item(132, 13).Value = -3246.473599999999  # M132 LeveProfitNQ
$ws.Cells.Item(132, 14).Value = -16911.2855  # N132 LeveProfitHQ

# Row 136: Metal with Mettle
$ws.Cells.Item(136, 8).Value = 2886.2083  # H136 currentAveragePrice
$ws.Cells.Item(136, 9).Value = 2697.4167  # I136 currentAveragePriceNQ
$ws.Cells.Item(136, 10).Value = 3452.5833  # J136 currentAveragePriceHQ
$ws.Cells.Item(136, 11).Value = 8092.250100000001  # K136 LevePriceNQ
$ws.Cells.Item(136, 12).Value = 10357.7499  # L136 LevePriceHQ
$ws.Cells.Item(136, 13).Value = -5542.250100000001  # M136 LeveProfitNQ
$ws.Cells.Item(136, 14).Value = -15457.7499  # N136 LeveProfitHQ

$ws = $wb.Worksheets.Item("BSM")
# Row 82: Spirituality Inspector
$ws.Cells.Item(82, 8).Value = 38864.785  # H82 currentAveragePrice
$ws.Cells.Item(82, 9).Value = 3868.625  # I82 currentAveragePriceNQ
$ws.Cells.Item(82, 11).Value = 3868.625  # K82 LevePriceNQ
$ws.Cells.Item(82, 13).Value = -3485.625  # M82 LeveProfitNQ

# Row 85: The Clamor for Hammers (L)
$ws.Cells.Item(85, 8).Value = 38864.785  # H85 currentAveragePrice
$ws.Cells.Item(85, 9).Value = 3868.625  # I85 currentAveragePriceNQ
$ws.Cells.Item(85, 11).Value = 3868.625  # K85 LevePriceNQ
$ws.Cells.Item(85, 13).Value = -2542.625  # M85 LeveProfitNQ

# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 2329.3914  # H86 currentAveragePrice
$ws.Cells.Item(86, 9).Value = 1920.6111  # I86 currentAveragePriceNQ
$ws.Cells.Item(86, 11).Value = 1920.6111  # K86 LevePriceNQ
$ws.Cells.Item(86, 13).Value = -797.6111000000001  # M86 LeveProfitNQ

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 2329.3914  # H89 currentAveragePrice
$ws.Cells.Item(89, 9).Value = 1920.6111  # I89 currentAveragePriceNQ
$ws.Cells.Item(89, 11).Value = 9603.0555  # K89 LevePriceNQ
$ws.Cells.Item(89, 13).Value = -3987.0555  # M89 LeveProfitNQ

# Row 94: High Steal
$ws.Cells.Item(94, 8).Value = 1083.8  # H94 currentAveragePrice
$ws.Cells.Item(94, 9).Value = 387.36365  # I94 currentAveragePriceNQ
$ws.Cells.Item(94, 11).Value = 387.36365  # K94 LevePriceNQ
$ws.Cells.Item(94, 13).Value = 63.63634999999999  # M94 LeveProfitNQ

# Row 134: Ruthenium Supremium
$ws.Cells.Item(134, 8).Value = 2113.5  # H134 currentAveragePrice
$ws.Cells.Item(134, 9).Value = 1121.2051  # I134 currentAveragePriceNQ
$ws.Cells.Item(134, 10).Value = 7642  # J134 currentAveragePriceHQ
$ws.Cells.Item(134, 11).Value = 3363.615299999999  # K134 LevePriceNQ
$ws.Cells.Item(134, 12).Value = 22926  # L134 LevePriceHQ
$ws.Cells.Item(134, 13).Value = -828.6152999999995  # M134 LeveProfitNQ
$ws.Cells.Item(134, 14).Value = -27996  # N134 LeveProfitHQ

$ws = $wb.Worksheets.Item("CRP")
# Row 6: Got Your Back
$ws.Cells.Item(6, 8).Value = 1033.8334  # H6 currentAveragePrice
$ws.Cells.Item(6, 9).Value = 1100.6  # I6 currentAveragePriceNQ
$ws.Cells.Item(6, 11).Value = 1100.6  # K6 LevePriceNQ
$ws.Cells.Item(6, 13).Value = -987.5999999999999  # M6 LeveProfitNQ

# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 81910.17  # H31 currentAveragePrice
$ws.Cells.Item(31, 9).Value = 108747.53  # I31 currentAveragePriceNQ
$ws.Cells.Item(31, 10).Value = 8717.362999999999  # J31 currentAveragePriceHQ
$ws.Cells.Item(31, 11).Value = 108747.53  # K31 LevePriceNQ
$ws.Cells.Item(31, 12).Value = 8717.362999999999  # L31 LevePriceHQ
$ws.Cells.Item(31, 13).Value = -108452.53  # M31 LeveProfitNQ
$ws.Cells.Item(31, 14).Value = -9307.362999999999  # N31 LeveProfitHQ

# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 81910.17  # H34 currentAveragePrice
$ws.Cells.Item(34, 9).Value = 108747.53  # I34 currentAveragePriceNQ
$ws.Cells.Item(34, 10).Value = 8717.362999999999  # J34 currentAveragePriceHQ
$ws.Cells.Item(34, 11).Value = 108747.53  # K34 LevePriceNQ
$ws.Cells.Item(34, 12).Value = 8717.362999999999  # L34 LevePriceHQ
$ws.Cells.Item(34, 13).Value = -108545.53  # M34 LeveProfitNQ
$ws.Cells.Item(34, 14).Value = -9121.362999999999  # N34 LeveProfitHQ

# Row 43: The Long Lance of the Law
$ws.Cells.Item(43, 8).Value = 81451.664  # H43 currentAveragePrice
$ws.Cells.Item(43, 10).Value = 81451.664  # J43 currentAveragePriceHQ
$ws.Cells.Item(43, 12).Value = 81451.664  # L43 LevePriceHQ
$ws.Cells.Item(43, 14).Value = -81819.664  # N43 LeveProfitHQ

# Row 58: You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 1913.9767  # H58 currentAveragePrice
$ws.Cells.Item(58, 9).Value = 1293.6428  # I58 currentAveragePriceNQ
$ws.Cells.Item(58, 10).Value = 3071.9333  # J58 currentAveragePriceHQ
$ws.Cells.Item(58, 11).Value = 1293.6428  # K58 LevePriceNQ
$ws.Cells.Item(58, 12).Value = 3071.9333  # L58 LevePriceHQ
$ws.Cells.Item(58, 13).Value = -1090.6428  # M58 LeveProfitNQ
$ws.Cells.Item(58, 14).Value = -3477.9333  # N58 LeveProfitHQ

# Row 86: Birch, Please
$ws.Cells.Item(86, 8).Value = 4501.75  # H86 currentAveragePrice
$ws.Cells.Item(86, 9).Value = 3007  # I86 currentAveragePriceNQ
$ws.Cells.Item(86, 11).Value = 3007  # K86 LevePriceNQ
$ws.Cells.Item(86, 13).Value = -1884  # M86 LeveProfitNQ

# Row 89: Built This City on Blocks and Soul (L)
$ws.Cells.Item(89, 8).Value = 4501.75  # H89 currentAveragePrice
$ws.Cells.Item(89, 9).Value = 3007  # I89 currentAveragePriceNQ
$ws.Cells.Item(89, 11).Value = 15035  # K89 LevePriceNQ
$ws.Cells.Item(89, 13).Value = -9419  # M89 LeveProfitNQ

# Row 99: O Pine
$ws.Cells.Item(99, 8).Value = 2750  # H99 currentAveragePrice
$ws.Cells.Item(99, 9).Value = 1500  # I99 currentAveragePriceNQ
$ws.Cells.Item(99, 10).Value = 4000  # J99 currentAveragePriceHQ
$ws.Cells.Item(99, 11).Value = 1500  # K99 LevePriceNQ
$ws.Cells.Item(99, 12).Value = 4000  # L99 LevePriceHQ
$ws.Cells.Item(99, 13).Value = -2  # M99 LeveProfitNQ
$ws.Cells.Item(99, 14).Value = -6996  # N99 LeveProfitHQ

# Row 101: Everybody's Heard about the 'Berd
$ws.Cells.Item(101, 8).Value = 81451.664  # H101 currentAveragePrice
$ws.Cells.Item(101, 10).Value = 81451.664  # J101 currentAveragePriceHQ
$ws.Cells.Item(101, 12).Value = 81451.664  # L101 LevePriceHQ
$ws.Cells.Item(101, 14).Value = -87941.664  # N101 LeveProfitHQ

# Row 126: A Better Conductor
$ws.Cells.Item(126, 8).Value = 2750  # H126 currentAveragePrice
$ws.Cells.Item(126, 9).Value = 1500  # I126 currentAveragePriceNQ
$ws.Cells.Item(126, 10).Value = 4000  # J126 currentAveragePriceHQ
$ws.Cells.Item(126, 11).Value = 4500  # K126 LevePriceNQ
$ws.Cells.Item(126, 12).Value = 12000  # L126 LevePriceHQ
$ws.Cells.Item(126, 13).Value = -2030  # M126 LeveProfitNQ
$ws.Cells.Item(126, 14).Value = -16940  # N126 LeveProfitHQ

# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 11906516  # H132 currentAveragePrice
$ws.Cells.Item(132, 9).Value = 1695.8235  # I132 currentAveragePriceNQ
$ws.Cells.Item(132, 10).Value = 62502004  # J132 currentAveragePriceHQ
$ws.Cells.Item(132, 11).Value = 5087.470499999999  # K132 LevePriceNQ
$ws.Cells.Item(132, 12).Value = 187506012  # L132 LevePriceHQ
$ws.Cells.Item(132, 13).Value = -2557.470499999999  # M132 LeveProfitNQ
$ws.Cells.Item(132, 14).Value = -187511072  # N132 LeveProfitHQ

# Row 136: Turali Quality
$ws.Cells.Item(136, 8).Value = 1913.9767  # H136 currentAveragePrice
$ws.Cells.Item(136, 9).Value = 1293.6428  # I136 currentAveragePriceNQ
$ws.Cells.Item(136, 10).Value = 3071.9333  # J136 currentAveragePriceHQ
$ws.Cells.Item(136, 11).Value = 3880.9284  # K136 LevePriceNQ
$ws.Cells.Item(136, 12).Value = 9215.7999  # L136 LevePriceHQ
$ws.Cells.Item(136, 13).Value = -1330.9284  # M136 LeveProfitNQ
$ws.Cells.Item(136, 14).Value = -14315.7999  # N136 LeveProfitHQ

$ws = $wb.Worksheets.Item("CUL")
# Row 86: Let's Not Get Sappy
$ws.Cells.Item(86, 8).Value = 408.7143  # H86 currentAveragePrice
$ws.Cells.Item(86, 10).Value = 144.2  # J86 currentAveragePriceHQ
$ws.Cells.Item(86, 12).Value = 432.6  # L86 LevePriceHQ
$ws.Cells.Item(86, 14).Value = -2804.6  # N86 LeveProfitHQ

# Row 89: Luxury Spillover (L)
$ws.Cells.Item(89, 8).Value = 408.7143  # H89 currentAveragePrice
$ws.Cells.Item(89, 10).Value = 144.2  # J89 currentAveragePriceHQ
$ws.Cells.Item(89, 12).Value = 1297.8  # L89 LevePriceHQ
$ws.Cells.Item(89, 14).Value = -13153.8  # N89 LeveProfitHQ

# Row 97: The Frier Never Lies
$ws.Cells.Item(97, 8).Value = 588  # H97 currentAveragePrice
$ws.Cells.Item(97, 10).Value = 588  # J97 currentAveragePriceHQ
$ws.Cells.Item(97, 12).Value = 1764  # L97 LevePriceHQ
$ws.Cells.Item(97, 14).Value = -2756  # N97 LeveProfitHQ

# Row 105: Fish Box
$ws.Cells.Item(105, 8).Value = 17542.143  # H105 currentAveragePrice
$ws.Cells.Item(105, 10).Value = 17542.143  # J105 currentAveragePriceHQ
$ws.Cells.Item(105, 12).Value = 52626.429  # L105 LevePriceHQ
$ws.Cells.Item(105, 14).Value = -57868.429  # N105 LeveProfitHQ

# Row 131: The Mountain Steeped
$ws.Cells.Item(131, 8).Value = 1462.7115  # H131 currentAveragePrice
$ws.Cells.Item(131, 9).Value = 915  # I131 currentAveragePriceNQ
$ws.Cells.Item(131, 10).Value = 1484.62  # J131 currentAveragePriceHQ
$ws.Cells.Item(131, 11).Value = 2745  # K131 LevePriceNQ
$ws.Cells.Item(131, 12).Value = 4453.86  # L131 LevePriceHQ
$ws.Cells.Item(131, 13).Value = 2295  # M131 LeveProfitNQ
$ws.Cells.Item(131, 14).Value = -14533.86  # N131 LeveProfitHQ

# Row 141: Ocean Explosion
$ws.Cells.Item(141, 8).Value = 2631.1428  # H141 currentAveragePrice
$ws.Cells.Item(141, 9).Value = 2564.3076  # I141 currentAveragePriceNQ
$ws.Cells.Item(141, 11).Value = 7692.9228  # K141 LevePriceNQ
$ws.Cells.Item(141, 13).Value = -2512.9228  # M141 LeveProfitNQ

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Cells.Item(80, 8).Value = 300133.47  # H80 currentAveragePrice
$ws.Cells.Item(80, 9).Value = 421864.16  # I80 currentAveragePriceNQ
$ws.Cells.Item(80, 10).Value = 7979.8  # J80 currentAveragePriceHQ
$ws.Cells.Item(80, 11).Value = 421864.16  # K80 LevePriceNQ
$ws.Cells.Item(80, 12).Value = 7979.8  # L80 LevePriceHQ
$ws.Cells.Item(80, 13).Value = -420866.16  # M80 LeveProfitNQ
$ws.Cells.Item(80, 14).Value = -9975.799999999999  # N80 LeveProfitHQ

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Cells.Item(83, 8).Value = 300133.47  # H83 currentAveragePrice
$ws.Cells.Item(83, 9).Value = 421864.16  # I83 currentAveragePriceNQ
$ws.Cells.Item(83, 10).Value = 7979.8  # J83 currentAveragePriceHQ
$ws.Cells.Item(83, 11).Value = 2109320.8  # K83 LevePriceNQ
$ws.Cells.Item(83, 12).Value = 39899  # L83 LevePriceHQ
$ws.Cells.Item(83, 13).Value = -2104328.8  # M83 LeveProfitNQ
$ws.Cells.Item(83, 14).Value = -49883  # N83 LeveProfitHQ

# Row 93: One Ring Circus
$ws.Cells.Item(93, 8).Value = 49999  # H93 currentAveragePrice
$ws.Cells.Item(93, 10).Value = 49999  # J93 currentAveragePriceHQ
$ws.Cells.Item(93, 12).Value = 49999  # L93 LevePriceHQ
$ws.Cells.Item(93, 14).Value = -53743  # N93 LeveProfitHQ

# Row 102: Put the Metal to the Peddle
$ws.Cells.Item(102, 8).Value = 4199.364  # H102 currentAveragePrice
$ws.Cells.Item(102, 9).Value = 4170.4287  # I102 currentAveragePriceNQ
$ws.Cells.Item(102, 10).Value = 4250  # J102 currentAveragePriceHQ
$ws.Cells.Item(102, 11).Value = 4170.4287  # K102 LevePriceNQ
$ws.Cells.Item(102, 12).Value = 4250  # L102 LevePriceHQ
$ws.Cells.Item(102, 13).Value = -2548.4287  # M102 LeveProfitNQ
$ws.Cells.Item(102, 14).Value = -7494  # N102 LeveProfitHQ

# Row 132: On Board for Lar
$ws.Cells.Item(132, 8).Value = 22741022  # H132 currentAveragePrice
$ws.Cells.Item(132, 9).Value = 33344998  # I132 currentAveragePriceNQ
$ws.Cells.Item(132, 10).Value = 18219.072  # J132 currentAveragePriceHQ
$ws.Cells.Item(132, 11).Value = 100034994  # K132 LevePriceNQ
$ws.Cells.Item(132, 12).Value = 54657.216  # L132 LevePriceHQ
$ws.Cells.Item(132, 13).Value = -100032464  # M132 LeveProfitNQ
$ws.Cells.Item(132, 14).Value = -59717.216  # N132 LeveProfitHQ

$ws = $wb.Worksheets.Item("LTW")
# Row 12: A Place to Call Helm
$ws.Cells.Item(12, 8).Value = 1115.5  # H12 currentAveragePrice
$ws.Cells.Item(12, 9).Value = 1298.6  # I12 currentAveragePriceNQ
$ws.Cells.Item(12, 10).Value = 200  # J12 currentAveragePriceHQ
$ws.Cells.Item(12, 11).Value = 1298.6  # K12 LevePriceNQ
$ws.Cells.Item(12, 12).Value = 200  # L12 LevePriceHQ
$ws.Cells.Item(12, 13).Value = -1128.6  # M12 LeveProfitNQ
$ws.Cells.Item(12, 14).Value = -540  # N12 LeveProfitHQ

# Row 22: Skin off Their Backs
$ws.Cells.Item(22, 8).Value = 3400.5  # H22 currentAveragePrice
$ws.Cells.Item(22, 9).Value = 3400.5  # I22 currentAveragePriceNQ
$ws.Cells.Item(22, 11).Value = 3400.5  # K22 LevePriceNQ
$ws.Cells.Item(22, 13).Value = -3105.5  # M22 LeveProfitNQ

# Row 27: Fire and Hide
$ws.Cells.Item(27, 8).Value = 3400.5  # H27 currentAveragePrice
$ws.Cells.Item(27, 9).Value = 3400.5  # I27 currentAveragePriceNQ
$ws.Cells.Item(27, 11).Value = 3400.5  # K27 LevePriceNQ
$ws.Cells.Item(27, 13).Value = -3293.5  # M27 LeveProfitNQ

# Row 46: Supply Side Logic
$ws.Cells.Item(46, 8).Value = 2887.05  # H46 currentAveragePrice
$ws.Cells.Item(46, 9).Value = 3027.5  # I46 currentAveragePriceNQ
$ws.Cells.Item(46, 10).Value = 2851.9375  # J46 currentAveragePriceHQ
$ws.Cells.Item(46, 11).Value = 3027.5  # K46 LevePriceNQ
$ws.Cells.Item(46, 12).Value = 2851.9375  # L46 LevePriceHQ
$ws.Cells.Item(46, 13).Value = -2839.5  # M46 LeveProfitNQ
$ws.Cells.Item(46, 14).Value = -3227.9375  # N46 LeveProfitHQ

# Row 62: Pummeling Abroad
$ws.Cells.Item(62, 8).Value = 41263  # H62 currentAveragePrice
$ws.Cells.Item(62, 9).Value = 23789  # I62 currentAveragePriceNQ
$ws.Cells.Item(62, 10).Value = 50000  # J62 currentAveragePriceHQ
$ws.Cells.Item(62, 11).Value = 23789  # K62 LevePriceNQ
$ws.Cells.Item(62, 12).Value = 50000  # L62 LevePriceHQ
$ws.Cells.Item(62, 13).Value = -23165  # M62 LeveProfitNQ
$ws.Cells.Item(62, 14).Value = -51248  # N62 LeveProfitHQ

# Row 65: The Style of the Time (L)
$ws.Cells.Item(65, 8).Value = 41263  # H65 currentAveragePrice
$ws.Cells.Item(65, 9).Value = 23789  # I65 currentAveragePriceNQ
$ws.Cells.Item(65, 10).Value = 50000  # J65 currentAveragePriceHQ
$ws.Cells.Item(65, 11).Value = 71367  # K65 LevePriceNQ
$ws.Cells.Item(65, 12).Value = 150000  # L65 LevePriceHQ
$ws.Cells.Item(65, 13).Value = -68247  # M65 LeveProfitNQ
$ws.Cells.Item(65, 14).Value = -156240  # N65 LeveProfitHQ

# Row 103: Security Breeches
$ws.Cells.Item(103, 8).Value = 86000.2  # H103 currentAveragePrice
$ws.Cells.Item(103, 10).Value = 86000.2  # J103 currentAveragePriceHQ
$ws.Cells.Item(103, 12).Value = 86000.2  # L103 LevePriceHQ
$ws.Cells.Item(103, 14).Value = -88344.2  # N103 LeveProfitHQ

$ws = $wb.Worksheets.Item("WVR")
# Row 14: Hat in Hand
$ws.Cells.Item(14, 8).Value = 4613.7856  # H14 currentAveragePrice
$ws.Cells.Item(14, 9).Value = 4549.4165  # I14 currentAveragePriceNQ
$ws.Cells.Item(14, 11).Value = 4549.4165  # K14 LevePriceNQ
$ws.Cells.Item(14, 13).Value = -4381.4165  # M14 LeveProfitNQ

# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 2026.973  # H132 currentAveragePrice
$ws.Cells.Item(132, 9).Value = 1905.1  # I132 currentAveragePriceNQ
$ws.Cells.Item(132, 11).Value = 5715.299999999999  # K132 LevePriceNQ
$ws.Cells.Item(132, 13).Value = -3185.299999999999  # M132 LeveProfitNQ

# Row 136: Weaving the Envelope
$ws.Cells.Item(136, 8).Value = 9776.706  # H136 currentAveragePrice
$ws.Cells.Item(136, 9).Value = 3153.7073  # I136 currentAveragePriceNQ
$ws.Cells.Item(136, 11).Value = 9461.1219  # K136 LevePriceNQ
$ws.Cells.Item(136, 13).Value = -6911.1219  # M136 LeveProfitNQ
